# The table (originally header in row1 B:F, data rows in rows2-6 with an
# index column in A) is reorganized: the index column A is dropped, and the
# remaining tickers/eps/revenue/price/people block (B1:F6) is shifted two
# columns to the right and one row down, landing at C2:G7.
#
# Because source (B1:F6) and destination (C2:G7) overlap, we stage the
# original block in a scratch area far away first, wipe the sheet, then
# copy the staged data into its final location (preserving values AND
# formatting), and finally clean up the scratch area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stage the original header+data block (including formatting) far away.
$ws.Range("B1:F6").Copy($ws.Range("Z1"))

# 2) Clear the original occupied area (old index column + header + data).
$ws.Range("A1:F6").Clear()

# 3) Copy the staged block into its new home, two columns right, one row down.
$ws.Range("Z1:AD6").Copy($ws.Range("C2"))

# 4) Clean up the scratch area so it doesn't linger in the saved workbook.
$ws.Range("Z1:AD6").Clear()
